$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16: risk table note.
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "risk table"
$ws.Range("C16").Value = "sudung table risk"

# New row 17: issue table note.
$ws.Range("A17").Value = 12
$ws.Range("B17").Value = "issue table"
$ws.Range("C17").Value = "su dung table issue"
$ws.Range("D17").Value = "workUnitID will be projectID"

# Update the "Role in Project" note row (row 3): new wording + a new
# multi-line explanation of the role values.
$ws.Range("C3").Value = "Xai bang Assigment role la Type"
$ws.Range("D3").Value = "0: ProjectOwer + PM, 1 : pM, 2: dev, 3: test, 4 QA, 5 cus, 6: project Owner`nProject Owner chi co' quyen read only va change PM,`nUser tao project se~ mang role la 0;"

# Make room for the multi-line note and wrap the text so it is readable.
$ws.Rows.Item(3).RowHeight = 45
$ws.Range("D3").WrapText = $true

# Leave the selection on D3, matching the author's last-edited cell.
$ws.Range("D3").Select()
